# Insert a new data row before the current row 469, shifting the existing
# rows 469-498 down to 470-499, then populate the newly inserted row 469
# with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 469 (pushes old 469..498 down to 470..499).
$ws.Rows.Item(469).Insert()

# Fill the newly inserted row with the new record's data.
$ws.Cells.Item(469, 1).Value = 3
$ws.Cells.Item(469, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(469, 3).Value = "Coquimbo"
$ws.Cells.Item(469, 4).Value = 44931
$ws.Cells.Item(469, 5).Value = 5
$ws.Cells.Item(469, 6).Value = 100112031
$ws.Cells.Item(469, 7).Value = "Poroto verde"
$ws.Cells.Item(469, 8).Value = "Magnum"
$ws.Cells.Item(469, 9).Value = "Primera"
$ws.Cells.Item(469, 10).Value = 78
$ws.Cells.Item(469, 11).Value = 30000
$ws.Cells.Item(469, 12).Value = 31000
$ws.Cells.Item(469, 13).Value = 30513
$ws.Cells.Item(469, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(469, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(469, 16).Value = 1221
$ws.Cells.Item(469, 17).Value = 25
$ws.Cells.Item(469, 18).Value = "Hortaliza"
